$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.67"
$ws.Range("E2").Value = "'0.78%"
$ws.Range("G2").Value = "'16"
$ws.Range("D3").Value = "'26.79"
$ws.Range("E3").Value = "'-1.12%"
$ws.Range("G3").Value = "'16"
$ws.Range("D4").Value = "'4.693"
$ws.Range("E4").Value = "'1.54%"
$ws.Range("G4").Value = "'16"
$ws.Range("D5").Value = "'0.05951"
$ws.Range("E5").Value = "'0.97%"
$ws.Range("G5").Value = "'16"
$ws.Range("D6").Value = "'6.635"
$ws.Range("E6").Value = "'-0.14%"
$ws.Range("G6").Value = "'16"
$ws.Range("D7").Value = "'0.8585"
$ws.Range("E7").Value = "'-0.67%"
$ws.Range("G7").Value = "'16"
$ws.Range("D8").Value = "'0.9247"
$ws.Range("E8").Value = "'-0.68%"
$ws.Range("G8").Value = "'16"
$ws.Range("E9").Value = "'-0.95%"
$ws.Range("G9").Value = "'16"
$ws.Range("D10").Value = "'0.04350"
$ws.Range("E10").Value = "'14.05%"
$ws.Range("G10").Value = "'16"
$ws.Range("D11").Value = "'0.07014"
$ws.Range("E11").Value = "'-1.11%"
$ws.Range("G11").Value = "'16"
$ws.Range("D12").Value = "'0.03007"
$ws.Range("E12").Value = "'-6.05%"
$ws.Range("G12").Value = "'16"
$ws.Range("D13").Value = "'0.09124"
$ws.Range("E13").Value = "'-1.08%"
$ws.Range("G13").Value = "'16"
$ws.Range("D14").Value = "'0.001528"
$ws.Range("E14").Value = "'-1.06%"
$ws.Range("G14").Value = "'16"
$ws.Range("D15").Value = "'0.0006050"
$ws.Range("E15").Value = "'-94.18%"
$ws.Range("G15").Value = "'16"
$ws.Range("D16").Value = "'0.006039"
$ws.Range("E16").Value = "'-1.27%"
$ws.Range("G16").Value = "'16"
$ws.Range("E17").Value = "'-1.24%"
$ws.Range("G17").Value = "'16"
$ws.Range("E18").Value = "'-1.91%"
$ws.Range("G18").Value = "'16"
$ws.Range("D19").Value = "'2.154"
$ws.Range("E19").Value = "'-2.62%"
$ws.Range("G19").Value = "'16"
$ws.Range("D20").Value = "'0.3101"
$ws.Range("E20").Value = "'-0.50%"
$ws.Range("G20").Value = "'16"
$ws.Range("D21").Value = "'0.1296"
$ws.Range("E21").Value = "'1.61%"
$ws.Range("G21").Value = "'16"
$ws.Range("D22").Value = "'3.812"
$ws.Range("E22").Value = "'-0.94%"
$ws.Range("G22").Value = "'16"
$ws.Range("D23").Value = "'0.04201"
$ws.Range("E23").Value = "'-0.21%"
$ws.Range("G23").Value = "'16"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'-0.49%"
$ws.Range("G24").Value = "'16"
$ws.Range("D25").Value = "'0.004481"
$ws.Range("E25").Value = "'4.66%"
$ws.Range("G25").Value = "'16"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("G26").Value = "'16"
$ws.Range("E27").Value = "'-11.58%"
$ws.Range("G27").Value = "'16"
$ws.Range("G28").Value = "'16"
$ws.Range("G29").Value = "'16"
$ws.Range("G30").Value = "'16"
$ws.Range("G31").Value = "'16"
$ws.Range("G32").Value = "'16"
$ws.Range("G33").Value = "'16"
$ws.Range("G34").Value = "'16"
$ws.Range("G35").Value = "'16"
$ws.Range("G36").Value = "'16"
$ws.Range("G37").Value = "'16"
$ws.Range("G38").Value = "'16"
$ws.Range("G39").Value = "'16"
$ws.Range("D40").Value = "'0.03818"
$ws.Range("E40").Value = "'-0.21%"
$ws.Range("G40").Value = "'16"
$ws.Range("D41").Value = "'0.1109"
$ws.Range("E41").Value = "'0.77%"
$ws.Range("G41").Value = "'16"
$ws.Range("D42").Value = "'0.003769"
$ws.Range("E42").Value = "'-38.08%"
$ws.Range("G42").Value = "'16"
$ws.Range("D43").Value = "'0.002427"
$ws.Range("E43").Value = "'8.03%"
$ws.Range("G43").Value = "'16"
$ws.Range("D44").Value = "'0.01493"
$ws.Range("E44").Value = "'31.82%"
$ws.Range("G44").Value = "'16"
$ws.Range("D45").Value = "'0.00005143"
$ws.Range("E45").Value = "'-5.71%"
$ws.Range("G45").Value = "'16"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("G46").Value = "'16"
$ws.Range("D47").Value = "'0.04994"
$ws.Range("E47").Value = "'-17.00%"
$ws.Range("G47").Value = "'16"
$ws.Range("D48").Value = "'0.2380"
$ws.Range("E48").Value = "'10,346.03%"
$ws.Range("G48").Value = "'16"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("G49").Value = "'16"
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("G50").Value = "'16"
$ws.Range("G51").Value = "'16"
